$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text (matching the source inlineStr cells) by
# temporarily applying a text number-format, then resetting the style so no
# extra formatting is left behind on the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D5").Value = "560.82"
$ws.Range("D6").Value = "137.27"
$ws.Range("D11").Value = "5.31"
$ws.Range("D14").Value = "33.55"
$ws.Range("D17").Value = "7.12"
$ws.Range("D20").Value = "432.65"
$ws.Range("D21").Value = "13.50"
$ws.Range("D22").Value = "0.716"
$ws.Range("D23").Value = "13.27"
$ws.Range("D25").Value = "79.61"
$ws.Range("D26").Value = "0.999"
$ws.Range("D30").Value = "7.70"
$ws.Range("D32").Value = "6.21"
$ws.Range("D33").Value = "25.61"
$ws.Range("D34").Value = "0.990"
$ws.Range("D35").Value = "5.86"
$ws.Range("D38").Value = "48.52"
$ws.Range("D39").Value = "8.64"
$ws.Range("D40").Value = "2.75"
$ws.Range("D41").Value = "396.81"
$ws.Range("D42").Value = "0.0351"
$ws.Range("D44").Value = "0.104"
$ws.Range("D47").Value = "122.58"
$ws.Range("D48").Value = "34.20"
$ws.Range("D51").Value = "23.23"

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"

# Remaining cells: safe to assign directly (non-numeric-looking text).
$ws.Range("D2").Value = "59.160.29"
$ws.Range("E2").Value = "  +2.41%  "
$ws.Range("D3").Value = "2.965.86"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("E6").Value = "  +4.11%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").Value = "2.959.29"
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("E10").Value = "  +3.99%  "
$ws.Range("E11").Value = "  +11.39%  "
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("E13").Value = "  +3.70%  "
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").Value = "3.456.04"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("D18").Value = "2.967.87"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("D19").Value = "59.275.39"
$ws.Range("E19").Value = "  +2.62%  "
$ws.Range("E20").Value = "  +4.28%  "
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("E22").Value = "  +3.25%  "
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("E27").Value = "  +10.36%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("E30").Value = "  +4.56%  "
$ws.Range("E31").Value = "  +8.31%  "
$ws.Range("E32").Value = "  +5.02%  "
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("E34").Value = "  +5.92%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E35").Value = "  +3.86%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0760"
$ws.Range("E36").Value = "  +9.84%  "
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("E40").Value = "  +5.65%  "
$ws.Range("E41").Value = "  +6.06%  "
$ws.Range("E42").Value = "  +2.03%  "
$ws.Range("D43").Value = "2.744.86"
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("E44").Value = "  -2.83%  "
$ws.Range("E45").Value = "  +6.03%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("E48").Value = "  +17.82%  "
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("E50").Value = "  +2.46%  "
$ws.Range("E51").Value = "  +1.67%  "
